$d = $word.ActiveDocument

# Locate the unique code line "const itens = []" (the variable declaration),
# so we don't touch the many unrelated Portuguese comment lines that also
# contain the word "itens".
$line = $d.Content
$found = $line.Find.Execute("const itens = []")
if (-not $found) {
    Write-Host "ERROR: could not find target line"
}
$lineStart = $line.Start

# "itens" sits at a fixed offset inside that literal line.
$target = $d.Range($lineStart + 6, $lineStart + 11)

# Rename itens -> listItens, keeping the run's existing character formatting
# (font, color, size, language) untouched.
$target.Text = "listItens"

# The renamed run now reads "listItens" (length 9) starting where "itens" used
# to start. Word's live spell-checker treats "listItens" as a possible
# misspelling/camel-case word and wraps it with proofErr spellStart/spellEnd
# markers, which in practice causes the run to be split into two pieces at
# the point the checker re-evaluated the text ("listI" + "tens").
$renamed = $d.Range($lineStart + 6, $lineStart + 15)
Write-Host "Renamed range text: $($renamed.Text)"

$secondPart = $d.Range($lineStart + 11, $lineStart + 15)
Write-Host "Second part text: $($secondPart.Text)"

# Force the paragraph's run list to split "listItens" into two runs right at
# the "listI" / "tens" boundary, without altering the resulting visible
# formatting: briefly nudge the font color of the trailing part away from its
# current value and immediately restore it. The restore re-applies the exact
# original color, so the final formatting for both runs is identical (same
# as the original "itens" run), while the run boundary remains in place.
$originalColor = $secondPart.Font.Color
$secondPart.Font.Color = 1
$secondPart.Font.Color = $originalColor

Write-Host "Final paragraph text: $($d.Range($lineStart, $lineStart + 15).Text)"
